$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NPLQ fixtures")

# Update the semi-final (SF) scores that were previously placeholders (-1)
$ws.Range("E131").Value = 1
$ws.Range("H131").Value = 4

$ws.Range("E132").Value = 0
$ws.Range("H132").Value = 0

$ws.Range("E133").Value = 2
$ws.Range("H133").Value = 1

$ws.Range("E134").Value = 2
$ws.Range("H134").Value = 3

$ws.Range("E135").Value = 3
$ws.Range("H135").Value = 1

# Row 136 (GF) previously had the wrong teams copied from the SF row above it.
# Correct it to be Gold Coast Knights vs Olympic FC.
$ws.Range("C136").Value = "Gold Coast Knights"
$ws.Range("D136").Value = "GCK"
$ws.Range("F136").Value = "Olympic FC"
$ws.Range("G136").Value = "BOL"

# Update the view to reflect scrolling down to the newly-updated rows
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 109
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H133").Select()
